$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the signature block (rows 28-29) down to rows 36-37 to make room
#    for the expanded worker/period table (grows from 8 to 16 data rows).
# ---------------------------------------------------------------------------
$ws.Range("B28:C28").UnMerge()
$ws.Range("B29:C29").UnMerge()
$ws.Range("H28:J28").UnMerge()
$ws.Range("H29:J29").UnMerge()

$ws.Range("B28:C29").Cut($ws.Range("B36"))
$ws.Range("H28:J29").Cut($ws.Range("H36"))

$ws.Rows("28:35").ClearContents()
$ws.Rows("28:35").ClearFormats()

$ws.Range("B36:C36").Merge()
$ws.Range("B37:C37").Merge()
$ws.Range("H36:J36").Merge()
$ws.Range("H37:J37").Merge()

# ---------------------------------------------------------------------------
# 2. Re-apply row formatting for the detail table.
#    Row 23 currently carries the "closing" border style used for the last
#    row of the table; push that style down to the new last row (31) first,
#    then stamp the regular row style (from row 16) across rows 17-30.
# ---------------------------------------------------------------------------
$ws.Range("B23:J23").Copy()
$ws.Range("B31:J31").PasteSpecial(-4122)

$ws.Range("B16:J16").Copy()
$ws.Range("B17:J30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Write the new worker / period detail rows (16-31).
#    Four workers, each with four pending periods (2507, 2506, 2505, 2504).
# ---------------------------------------------------------------------------
$workers = @(
  @{ Doc = "73140520";   Name = "OSWALDO DE JESUS CASTILLA TARRA"; Salario = 1423500 },
  @{ Doc = "1102868229"; Name = "FREDY DE JESUS MENDOZA PEREZ";    Salario = 1423500 },
  @{ Doc = "1052079546"; Name = "MARIA JOSE GONZALEZ ANGULO";      Salario = 1423500 },
  @{ Doc = "45649374";   Name = "LICED MARGARITA TAPIA TORRES";    Salario = 1900000 }
)
$periodos = @("2507", "2506", "2505", "2504")

$row = 16
foreach ($w in $workers) {
  foreach ($p in $periodos) {
    $ws.Cells.Item($row, 2).Value() = "CC"
    $ws.Cells.Item($row, 3).Value() = $w.Doc
    $ws.Cells.Item($row, 4).Value() = $w.Name
    $ws.Cells.Item($row, 5).Value() = $p
    $ws.Cells.Item($row, 6).Value() = 56940
    $ws.Cells.Item($row, 7).Value() = $w.Salario
    $row = $row + 1
  }
}

# ---------------------------------------------------------------------------
# 4. Update the summary figures.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value() = 911040
$ws.Range("F13").Value() = 4
